$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '31.110.16'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.97%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.980.85'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.13%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.94%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.01%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7359'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +13.59%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.006'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.91%  '

# Row 8
$ws.Range("E8").Value = '  +4.22%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.47'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.42%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07122'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.79%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8281'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.54%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08106'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.65%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.990.66'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.73%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.586'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.96%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '99.29'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.67%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +11.57%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '268.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.39%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '31.129.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.99%  '

# Row 19
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008233'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.26%  '

# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.038'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.31%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.254.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.81%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.007'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.78%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.006'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.94%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.047'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.14%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.958'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.83%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.03'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.95%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.17%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.352'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1318'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.85%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.598'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.73%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.382'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.58%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.613'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.40%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.391'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.45%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05298'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.21%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.280'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.88%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7786'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.56%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.795'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.38%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.885'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.75%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '83.33'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.49%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.760'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.79%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4625'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.66%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.107'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.81%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8544'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.92%  '

# Row 45
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.74'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.12%  '

# Row 46
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.006'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.90%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.48%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.628'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.95%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.579'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +11.57%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.57%  '

# Row 51
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.949'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +38.37%  '
